# Final update to go client
# Re-sequence the Item Name (D) and UOM (E) columns on Sheet1. The BSL/BRAND/ISL
# columns (A,B,C) are unaffected; only the Item Name / UOM pairing per row changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @{ Row = 2;  D = "Desodin 60ml Syrup";                 E = "60 ml" },
    @{ Row = 3;  D = "Dinafex 180mg Tablet";                E = "30's" },
    @{ Row = 4;  D = "Dinafex 120mg Tablet";                E = "30's" },
    @{ Row = 5;  D = "Dinafex 60mg Tablet";                 E = "30's" },
    @{ Row = 6;  D = "Dorenta 50mg Tablet";                 E = "50's" },
    @{ Row = 7;  D = "Etorix 120mg Tablet";                 E = "20's" },
    @{ Row = 8;  D = "Etorix 60mg Tablet - 40's";           E = "40's" },
    @{ Row = 9;  D = "Etorix 90mg Tablet";                  E = "30's" },
    @{ Row = 10; D = "Fenobac 100ml Syrup";                 E = "100ml" },
    @{ Row = 11; D = "Flucloxin 500mg Capsule";             E = "30 's" },
    @{ Row = 12; D = "Flucloxin 500mg Capsule - 36's";      E = "36 's" },
    @{ Row = 13; D = "Geminox 320mg Tablet - 8's";          E = "8 's" },
    @{ Row = 14; D = "Ketonic 30mg IM/IV Injection - 4's";  E = "4's" },
    @{ Row = 15; D = "Ketonic 10mg Tablet";                 E = "20's" },
    @{ Row = 16; D = "Ketonic 30mg Injection";              E = "5 's" },
    @{ Row = 17; D = "Kynol TR 200mg Capsule";              E = "30 's" },
    @{ Row = 18; D = "Kynol TR 100mg Capsule";              E = "50 's" },
    @{ Row = 19; D = "Kynol D 25mg Tablet";                 E = "60 's" },
    @{ Row = 20; D = "Naprox Plus 500mg Tablet - 30's";     E = "30 's" },
    @{ Row = 21; D = "Oradin Plus Tablet - 40's";           E = "40 's" },
    @{ Row = 22; D = "Osticare Tablet 24's";                E = "24's" },
    @{ Row = 23; D = "Sk-Mox 500mg Capsule";                E = "48 's" },
    @{ Row = 24; D = "Zithrox 15ml Suspension";             E = "15 ml" },
    @{ Row = 25; D = "Zithrox 30ml Dry Suspension";         E = "30ml" },
    @{ Row = 26; D = "Zithrox 500mg Tablet";                E = "6 's" },
    @{ Row = 27; D = "Zithrox 250mg Tablet - 6's";          E = "6's" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
